# Horwath_2008.xlsx - "Updated soil type data"
#
# Insert a new "pro_usda_soil_order" column into the "profile" sheet
# (immediately after the "pro_MAP" column, i.e. before the previous
# "pro_soil_taxon" column), add its header, and record the USDA soil
# order ("Gelisols") for the existing data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("profile")

# Column N (14) is currently "pro_soil_taxon"; inserting here shifts
# pro_soil_taxon and everything after it one column to the right,
# opening up a blank column N for the new field.
$ws.Columns.Item(14).Insert()

# Row 1 holds the machine-readable column/variable names.
$ws.Cells.Item(1, 14).Value = "pro_usda_soil_order"

# Row 4 holds the Horwath_2008 data values; record the USDA soil order.
$ws.Cells.Item(4, 14).Value = "Gelisols"
